$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7466.1113
$ws.Range("I40").Value = 11690.1
$ws.Range("J40").Value = 2186.125
$ws.Range("K40").Value = 11690.1
$ws.Range("L40").Value = 2186.125
$ws.Range("M40").Value = -11515.1
$ws.Range("N40").Value = -2536.125
$ws.Range("H41").Value = 1216.6666
$ws.Range("I41").Value = 950
$ws.Range("J41").Value = 1350
$ws.Range("K41").Value = 950
$ws.Range("L41").Value = 1350
$ws.Range("M41").Value = -510
$ws.Range("N41").Value = -2230
$ws.Range("H86").Value = 4270.4707
$ws.Range("I86").Value = 4224.875
$ws.Range("J86").Value = 4311
$ws.Range("K86").Value = 4224.875
$ws.Range("L86").Value = 4311
$ws.Range("M86").Value = -3101.875
$ws.Range("N86").Value = -6557
$ws.Range("H89").Value = 4270.4707
$ws.Range("I89").Value = 4224.875
$ws.Range("J89").Value = 4311
$ws.Range("K89").Value = 21124.375
$ws.Range("L89").Value = 21555
$ws.Range("M89").Value = -15508.375
$ws.Range("N89").Value = -32787
$ws.Range("H114").Value = 43722
$ws.Range("J114").Value = 43722
$ws.Range("L114").Value = 43722
$ws.Range("N114").Value = -52400
$ws.Range("H126").Value = 45780
$ws.Range("J126").Value = 45780
$ws.Range("L126").Value = 45780
$ws.Range("N126").Value = -55660
$ws.Range("H129").Value = 994.94116
$ws.Range("J129").Value = 933.2763
$ws.Range("L129").Value = 2799.8289
$ws.Range("N129").Value = -12799.8289
$ws.Range("H137").Value = 3432.6365
$ws.Range("I137").Value = 972.9583
$ws.Range("K137").Value = 2918.8749
$ws.Range("M137").Value = -368.8748999999998

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1437.4286
$ws.Range("I61").Value = 1139.5
$ws.Range("J61").Value = 1973.7
$ws.Range("K61").Value = 1139.5
$ws.Range("L61").Value = 1973.7
$ws.Range("M61").Value = -927.5
$ws.Range("N61").Value = -2397.7
$ws.Range("H74").Value = 1684.5264
$ws.Range("I74").Value = 1411.25
$ws.Range("J74").Value = 2449.7
$ws.Range("K74").Value = 1411.25
$ws.Range("L74").Value = 2449.7
$ws.Range("M74").Value = -537.25
$ws.Range("N74").Value = -4197.7
$ws.Range("H77").Value = 1684.5264
$ws.Range("I77").Value = 1411.25
$ws.Range("J77").Value = 2449.7
$ws.Range("K77").Value = 7056.25
$ws.Range("L77").Value = 12248.5
$ws.Range("M77").Value = -2688.25
$ws.Range("N77").Value = -20984.5
$ws.Range("H132").Value = 15153761
$ws.Range("I132").Value = 22728344
$ws.Range("J132").Value = 4594.4546
$ws.Range("K132").Value = 68185032
$ws.Range("L132").Value = 13783.3638
$ws.Range("M132").Value = -68182502
$ws.Range("N132").Value = -18843.3638
$ws.Range("H136").Value = 1437.4286
$ws.Range("I136").Value = 1139.5
$ws.Range("J136").Value = 1973.7
$ws.Range("K136").Value = 3418.5
$ws.Range("L136").Value = 5921.1
$ws.Range("M136").Value = -868.5
$ws.Range("N136").Value = -11021.1

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3258.3618
$ws.Range("I134").Value = 1568
$ws.Range("J134").Value = 3658.7104
$ws.Range("K134").Value = 4704
$ws.Range("L134").Value = 10976.1312
$ws.Range("M134").Value = -2169
$ws.Range("N134").Value = -16046.1312

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1785.16
$ws.Range("I31").Value = 700.8461
$ws.Range("J31").Value = 2959.8333
$ws.Range("K31").Value = 700.8461
$ws.Range("L31").Value = 2959.8333
$ws.Range("M31").Value = -405.8461
$ws.Range("N31").Value = -3549.8333
$ws.Range("H34").Value = 1785.16
$ws.Range("I34").Value = 700.8461
$ws.Range("J34").Value = 2959.8333
$ws.Range("K34").Value = 700.8461
$ws.Range("L34").Value = 2959.8333
$ws.Range("M34").Value = -498.8461
$ws.Range("N34").Value = -3363.8333
$ws.Range("H58").Value = 1523.381
$ws.Range("I58").Value = 1318.1333
$ws.Range("J58").Value = 2036.5
$ws.Range("K58").Value = 1318.1333
$ws.Range("L58").Value = 2036.5
$ws.Range("M58").Value = -1115.1333
$ws.Range("N58").Value = -2442.5
$ws.Range("H132").Value = 142828.5
$ws.Range("I132").Value = 1047.5
$ws.Range("J132").Value = 237349.17
$ws.Range("K132").Value = 3142.5
$ws.Range("L132").Value = 712047.51
$ws.Range("M132").Value = -612.5
$ws.Range("N132").Value = -717107.51
$ws.Range("H134").Value = 540270.5600000001
$ws.Range("I134").Value = 1342.5883
$ws.Range("J134").Value = 1558245.6
$ws.Range("K134").Value = 4027.7649
$ws.Range("L134").Value = 4674736.800000001
$ws.Range("M134").Value = -1492.7649
$ws.Range("N134").Value = -4679806.800000001
$ws.Range("H136").Value = 1523.381
$ws.Range("I136").Value = 1318.1333
$ws.Range("J136").Value = 2036.5
$ws.Range("K136").Value = 3954.3999
$ws.Range("L136").Value = 6109.5
$ws.Range("M136").Value = -1404.3999
$ws.Range("N136").Value = -11209.5

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4503.7407
$ws.Range("J5").Value = 1640.7142
$ws.Range("L5").Value = 4922.142599999999
$ws.Range("N5").Value = -5146.142599999999
$ws.Range("H33").Value = 17796832
$ws.Range("I33").Value = 149.75
$ws.Range("J33").Value = 24268352
$ws.Range("K33").Value = 898.5
$ws.Range("L33").Value = 145610112
$ws.Range("M33").Value = -615.5
$ws.Range("N33").Value = -145610678
$ws.Range("H113").Value = 2354.9648
$ws.Range("J113").Value = 743.55554
$ws.Range("L113").Value = 2230.66662
$ws.Range("N113").Value = -6570.66662
$ws.Range("H135").Value = 4503.7407
$ws.Range("J135").Value = 1640.7142
$ws.Range("L135").Value = 14766.4278
$ws.Range("N135").Value = -19836.4278

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2533.2104
$ws.Range("I132").Value = 1644.2273
$ws.Range("K132").Value = 4932.6819
$ws.Range("M132").Value = -2402.6819

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 20000
$ws.Range("J81").Value = 20000
$ws.Range("L81").Value = 20000
$ws.Range("N81").Value = -21996
$ws.Range("H84").Value = 20000
$ws.Range("J84").Value = 20000
$ws.Range("L84").Value = 60000
$ws.Range("N84").Value = -69984
$ws.Range("H132").Value = 2403
$ws.Range("I132").Value = 1766.0834
$ws.Range("K132").Value = 5298.2502
$ws.Range("M132").Value = -2768.2502
$ws.Range("H136").Value = 2136.5386
$ws.Range("I136").Value = 1687.0526
$ws.Range("K136").Value = 5061.1578
$ws.Range("M136").Value = -2511.1578

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 547
$ws.Range("I100").Value = 359
$ws.Range("J100").Value = 609.6667
$ws.Range("K100").Value = 718
$ws.Range("L100").Value = 1219.3334
$ws.Range("M100").Value = -177
$ws.Range("N100").Value = -2301.3334
$ws.Range("H132").Value = 2088.842
$ws.Range("I132").Value = 1715.0834
$ws.Range("J132").Value = 2729.5715
$ws.Range("K132").Value = 5145.2502
$ws.Range("L132").Value = 8188.7145
$ws.Range("M132").Value = -2615.2502
$ws.Range("N132").Value = -13248.7145
$ws.Range("H136").Value = 371307.03
$ws.Range("I136").Value = 435494.56
$ws.Range("J136").Value = 2228.75
$ws.Range("K136").Value = 1306483.68
$ws.Range("L136").Value = 6686.25
$ws.Range("M136").Value = -1303933.68
$ws.Range("N136").Value = -11786.25
